$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (ALC)
$ws.Range("H2").Value = 320.53845
$ws.Range("I2").Value = 304.2
$ws.Range("J2").Value = 375
$ws.Range("K2").Value = 304.2
$ws.Range("L2").Value = 375
$ws.Range("M2").Value = -191.2
$ws.Range("N2").Value = -601

# Row 12 (ALC)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()

# Row 18 (ALC)
$ws.Range("H18").Value = 8053.846
$ws.Range("I18").Value = 11411.223
$ws.Range("J18").Value = 499.75
$ws.Range("K18").Value = 11411.223
$ws.Range("L18").Value = 499.75
$ws.Range("M18").Value = -11127.223
$ws.Range("N18").Value = -1067.75

# Row 28 (ALC)
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()

# Row 32 (ALC)
$ws.Range("H32").Value = 11000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 11000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 11000
$ws.Range("N32").Value = -11652
$ws.Range("M32").ClearContents()

# Row 61 (ALC)
$ws.Range("H61").Value = 4234
$ws.Range("I61").Value = 3973.3333
$ws.Range("J61").Value = 5016
$ws.Range("K61").Value = 11919.9999
$ws.Range("L61").Value = 15048
$ws.Range("M61").Value = -11747.9999
$ws.Range("N61").Value = -15392

# Row 98 (ALC)
$ws.Range("H98").Value = 3152.4
$ws.Range("I98").Value = 2280.4443
$ws.Range("J98").Value = 11000
$ws.Range("K98").Value = 2280.4443
$ws.Range("L98").Value = 11000
$ws.Range("M98").Value = -782.4443000000001
$ws.Range("N98").Value = -13996

# Row 122 (ALC)
$ws.Range("H122").Value = 3152.4
$ws.Range("I122").Value = 2280.4443
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 6841.3329
$ws.Range("L122").Value = 33000
$ws.Range("M122").Value = -4391.3329
$ws.Range("N122").Value = -37900

# Row 135 (ALC)
$ws.Range("H135").Value = 1330.75
$ws.Range("I135").Value = 774.5
$ws.Range("J135").Value = 2999.5
$ws.Range("K135").Value = 6970.5
$ws.Range("L135").Value = 26995.5
$ws.Range("M135").Value = -4435.5
$ws.Range("N135").Value = -32065.5

$ws = $wb.Worksheets.Item("ARM")
# Row 10 (ARM)
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1830

# Row 32 (ARM)
$ws.Range("H32").Value = 3206.913
$ws.Range("I32").Value = 3248.2273
$ws.Range("J32").Value = 2298
$ws.Range("K32").Value = 3248.2273
$ws.Range("L32").Value = 2298
$ws.Range("M32").Value = -2961.2273
$ws.Range("N32").Value = -2872

# Row 62 (ARM)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65 (ARM)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 97 (ARM)
$ws.Range("H97").Value = 4326.5713
$ws.Range("I97").Value = 4121.5
$ws.Range("J97").Value = 4600
$ws.Range("K97").Value = 4121.5
$ws.Range("L97").Value = 4600
$ws.Range("M97").Value = -3625.5
$ws.Range("N97").Value = -5592

$ws = $wb.Worksheets.Item("BSM")
# Row 132 (BSM)
$ws.Range("H132").Value = 99780
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 99780
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 99780
$ws.Range("N132").Value = -109900

$ws = $wb.Worksheets.Item("CRP")
# Row 3 (CRP)
$ws.Range("H3").Value = 2375
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 1750
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 1750
$ws.Range("M3").Value = -2887
$ws.Range("N3").Value = -1976

# Row 88 (CRP)
$ws.Range("H88").Value = 12548.4
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 12548.4
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 12548.4
$ws.Range("N88").Value = -13360.4

# Row 91 (CRP)
$ws.Range("H91").Value = 12548.4
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 12548.4
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 12548.4
$ws.Range("N91").Value = -15356.4

# Row 107 (CRP)
$ws.Range("H107").Value = 645.5
$ws.Range("I107").Value = 615.4
$ws.Range("J107").Value = 796
$ws.Range("K107").Value = 615.4
$ws.Range("L107").Value = 796
$ws.Range("M107").Value = 1304.6
$ws.Range("N107").Value = -4636

$ws = $wb.Worksheets.Item("CUL")
# Row 57 (CUL)
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 7 (GSM)
$ws.Range("H7").Value = 253200.5
$ws.Range("I7").Value = 501201
$ws.Range("J7").Value = 5200
$ws.Range("K7").Value = 501201
$ws.Range("L7").Value = 5200
$ws.Range("M7").Value = -501089
$ws.Range("N7").Value = -5424

# Row 8 (GSM)
$ws.Range("H8").Value = 253200.5
$ws.Range("I8").Value = 501201
$ws.Range("J8").Value = 5200
$ws.Range("K8").Value = 501201
$ws.Range("L8").Value = 5200
$ws.Range("M8").Value = -501062
$ws.Range("N8").Value = -5478

# Row 80 (GSM)
$ws.Range("H80").Value = 2585.1667
$ws.Range("I80").Value = 2585.1667
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2585.1667
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1587.1667

# Row 83 (GSM)
$ws.Range("H83").Value = 2585.1667
$ws.Range("I83").Value = 2585.1667
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12925.8335
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7933.833500000001

# Row 97 (GSM)
$ws.Range("H97").Value = 1496.125
$ws.Range("I97").Value = 529
$ws.Range("J97").Value = 3108
$ws.Range("K97").Value = 529
$ws.Range("L97").Value = 3108
$ws.Range("M97").Value = -33
$ws.Range("N97").Value = -4100

# Row 102 (GSM)
$ws.Range("H102").Value = 1831
$ws.Range("I102").Value = 1831
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1831
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -209

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 4418.1665
$ws.Range("I46").Value = 2750
$ws.Range("J46").Value = 4751.8
$ws.Range("K46").Value = 2750
$ws.Range("L46").Value = 4751.8
$ws.Range("M46").Value = -2562
$ws.Range("N46").Value = -5127.8

# Row 53 (LTW)
$ws.Range("H53").Value = 2000
$ws.Range("I53").Value = 2000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 2000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -1482
$ws.Range("N53").ClearContents()

# Row 62 (LTW)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 64 (LTW)
$ws.Range("H64").Value = 13999.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 13999.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 13999.5
$ws.Range("N64").Value = -14449.5

# Row 65 (LTW)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 67 (LTW)
$ws.Range("H67").Value = 13999.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 13999.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 13999.5
$ws.Range("N67").Value = -15559.5

# Row 76 (LTW)
$ws.Range("H76").Value = 19743
$ws.Range("I76").Value = 19286
$ws.Range("J76").Value = 20200
$ws.Range("K76").Value = 19286
$ws.Range("L76").Value = 20200
$ws.Range("M76").Value = -18948
$ws.Range("N76").Value = -20876

# Row 79 (LTW)
$ws.Range("H79").Value = 19743
$ws.Range("I79").Value = 19286
$ws.Range("J79").Value = 20200
$ws.Range("K79").Value = 19286
$ws.Range("L79").Value = 20200
$ws.Range("M79").Value = -18116
$ws.Range("N79").Value = -22540

# Row 87 (LTW)
$ws.Range("H87").Value = 49999
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 49999
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 49999
$ws.Range("N87").Value = -52245

# Row 90 (LTW)
$ws.Range("H90").Value = 49999
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 49999
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 149997
$ws.Range("N90").Value = -161229

# Row 100 (LTW)
$ws.Range("H100").Value = 2824.5
$ws.Range("I100").Value = 2932.6667
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 2932.6667
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -2391.6667
$ws.Range("N100").Value = -3582

# Row 122 (LTW)
$ws.Range("H122").Value = 4149
$ws.Range("I122").Value = 4365.3335
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 13096.0005
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -10646.0005
$ws.Range("N122").Value = -15400

# Row 136 (LTW)
$ws.Range("H136").Value = 1825.5714
$ws.Range("I136").Value = 1546.5
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4639.5
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2089.5
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("WVR")
# Row 55 (WVR)
$ws.Range("H55").Value = 1500
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 1000
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = -723
$ws.Range("N55").Value = -2554

# Row 63 (WVR)
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66 (WVR)
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
